$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; this shifts the old rows 4-13 down to 5-14,
# carrying their existing content/formatting along unchanged.
$ws.Rows("4:4").Insert()

# The newly inserted row only ends up holding data in columns B and E, so
# drop the (empty, inherited) formatting in the other columns of that row.
$ws.Range("A4").Clear()
$ws.Range("C4").Clear()
$ws.Range("D4").Clear()

# Row 2: "Assistant/Associate Professor" end-date becomes a literal year (2019)
# instead of the shared "2017 - Present" text.
$ws.Range("B2").Value = 2019

# Row 3: Assistant Professor entry now covers 2017-2018 (text) instead of
# just the single year 2017, and its description becomes the
# "Quantitative Methods II" text (same as row 2's description).
$ws.Range("B3").Value = "2017-2018"
$ws.Range("E3").Value = "Quantitative Methods II (Psychology MSc)."

# New row 4 (inherits formatting from the row above via the insert): holds
# the split-out "2017" entry with the "Quantitative Methods I" description
# that used to live on row 3.
$ws.Range("B4").Value = 2017
$ws.Range("E4").Value = "Quantitative Methods I (Psychology MSc)."

# Update the selection to match the saved view state.
$ws.Range("C16").Select()
